# Automatische test-sync: 2025-06-19 17:57:30
# Append a new incoming-mail log row to the "Logs" sheet and bump the
# matching category tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New row goes right after the current last used row (row 30 -> row 31).
$newRow = 31

$logs.Cells.Item($newRow, 1).Value = "Klacht over levering"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$logs.Cells.Item($newRow, 4).Value = "Klacht"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 17:57:20"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# The conditional-formatting ranges on columns D and G are anchored to the
# sheet's used range (2:30) and need to grow along with the new row.
$dRules = $logs.Range("D2:D30").FormatConditions
for ($i = 1; $i -le $dRules.Count; $i++) {
    $dRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D31"))
}

$gRules = $logs.Range("G2:G30").FormatConditions
for ($i = 1; $i -le $gRules.Count; $i++) {
    $gRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G31"))
}

# Bump the "Klacht" tally on the Dashboard sheet (5 -> 6) using the
# current numeric value (read via Value2 since Value's getter is unreliable
# in this host for round-tripping numbers).
$klachtCell = $dashboard.Range("B4")
$klachtCell.Value = $klachtCell.Value2 + 1
